# Move E3hr prrc from identified missing to ignored list (#444).
#
# Row 622 of Sheet1 holds the "E3hr" / "prrc" entry (with its long
# description in column H and "Twan & Thomas" in column I). That entire
# row is being dropped from the "identified missing" list, so every row
# below it shifts up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row so everything beneath it shifts up by one.
$ws.Rows.Item(622).Delete()

# Leave the selection/scroll position where the edit happened, matching
# the workbook's view state after the change.
$excel.ActiveWindow.ScrollRow = 605
$ws.Range("A622").Select()
